$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: "Largest Rectangle in Histogram" gains STATUS / LAST SOLVED / NOTES
$ws.Range("G28").Value = "STRUGGLED"
$ws.Range("H28").Value = "13/06/2025"
$ws.Range("I28").Value = "Tough One!!!"

# Row 29: new entry - Binary Search (Leetcode 704, Easy, Neetcode 150)
$ws.Range("A29").Value = "Leetcode"
$ws.Range("B29").Value = 704
$ws.Range("C29").Value = "Binary Search"
$ws.Range("D29").Value = "Binary Search"
$ws.Range("E29").Value = "Easy"
$ws.Range("F29").Value = "Neetcode 150"
$ws.Range("H29").Value = "14/06/2025"

# Scroll the frozen view down so row 9 is the first visible row under the
# frozen header (matches the author scrolling down to the newly added rows).
$excel.ActiveWindow.ScrollRow = 9
[void]$ws.Range("C30").Select()
